# Generate Report for Handoff
# Update status/handoff rows for the "9f4ad892-dc37-4aec-b18d-4595c510be76.md" file
# across the Overview, zh-cn, and de-de sheets, reflecting a new handoff event.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to 9f4ad892-dc37-4aec-b18d-4595c510be76.md
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 corresponds to 9f4ad892-dc37-4aec-b18d-4595c510be76.md
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-02-22 09:22:37"

# de-de sheet: row 3 corresponds to 9f4ad892-dc37-4aec-b18d-4595c510be76.md
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-02-22 09:22:49"
